$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for rows 2-15, columns A (date serial), B, C, D, E, F
$data = @(
    ,@(2, 45980, 10350.396094346, 9822.94242402898, 16491.4, 5568.92997055721, -45.8136502255755)
    ,@(3, 45981, 10452.4755653533, 9776.35955251518, 9003.4, 6043.0651391464, 284.001028819233)
    ,@(4, 45982, 10429.4640575845, 9312.23245280107, 9003.4, 6048.78186898908, 264.900596741257)
    ,@(5, 45983, 4407.8584487922, 7135.9851811023, 9003.4, 5974.28038066799, 171.119398407096)
    ,@(6, 45984, 4649.26780577935, 7486.49656737432, 9003.4, 6447.91159295917, 205.458673347229)
    ,@(7, 45985, 10980.8696790314, 10114.663427364, 9003.4, 6409.2357407555, 313.354132004981)
    ,@(8, 45986, 10980.8696790314, 10054.5062684567, 9003.4, 6409.2357407555, 310.847583717174)
    ,@(9, 45987, 10980.8696790314, 10414.192954408, 9003.4, 6409.2357407555, 325.834528965148)
    ,@(10, 45988, 10980.8696790314, 10376.8181471995, 9003.4, 6409.2357407555, 324.277245331459)
    ,@(11, 45989, 10980.8696790314, 9629.0457859519, 9003.4, 6409.2357407555, 293.120063612809)
    ,@(12, 45990, 4442.42485201967, 6726.19293427106, 9003.4, 6033.01536760749, 156.492012578273)
    ,@(13, 45991, 4286.74551609542, 6538.57638271881, 9003.4, 6024.52678564567, 148.32096534852)
    ,@(14, 45992, 11778.1169212215, 10656.5494641501, 8664.26, 7488.33918927062, 395.026193892531)
    ,@(15, 45993, 11778.1169212215, 11009.5282379234, 8664.26, 7488.33918927062, 409.73364279975)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
    $ws.Cells.Item($r, 5).Value = $row[5]
    $ws.Cells.Item($r, 6).Value = $row[6]
}
